# Swap species-record data between row 25 <-> row 27, and row 28 <-> row 29
# (columns A, B, D, E, F, G, H, Q, R) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

function Swap-Rows($ws, $cols, $rowA, $rowB) {
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value2 = $valB
        $rangeB.Value2 = $valA
    }
}

Swap-Rows $ws $cols 25 27
Swap-Rows $ws $cols 28 29
